# Insert a new data row above row 78 (shifts existing rows 78..171 down to 79..172)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78; this shifts rows 78:171 down to 79:172
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record's data
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"

# Date column (style copies default general number; set underlying numeric serial value
# and apply the same date number format used by the other rows in column D)
$ws.Cells.Item(78, 4).Value = 44895
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = 100112032
$ws.Cells.Item(78, 7).Value = "Zapallo italiano"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 100
$ws.Cells.Item(78, 11).Value = 7000
$ws.Cells.Item(78, 12).Value = 7500
$ws.Cells.Item(78, 13).Value = 7250
$ws.Cells.Item(78, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(78, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(78, 16).Value = 145
$ws.Cells.Item(78, 17).Value = 50
$ws.Cells.Item(78, 18).Value = "Hortaliza"
